$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ErrorMSG")

# Row 108: error text switched to the "not found" message (endpoint name C108 unchanged)
$ws.Cells.Item(108, 2).Value = "A keresett adatok nem találhatóak"

# Rows 111-112: CAPTCHA check rows now point at CheckCaptcha (unchanged text, already there)
$ws.Cells.Item(111, 3).Value = "CheckCaptcha"
$ws.Cells.Item(112, 3).Value = "CheckCaptcha"

# Rows 118-120: GetMatchHistory (new endpoint) replacing GetUserStatistic block
$ws.Cells.Item(118, 3).Value = "GetMatchHistory"

$ws.Cells.Item(119, 2).Value = "Kapcsolati hiba"
$ws.Cells.Item(119, 3).Value = "GetMatchHistory"

$ws.Cells.Item(120, 2).Value = "Hiba történt az adatok mentése közben"
$ws.Cells.Item(120, 3).Value = "GetMatchHistory"

# Rows 121-123: GetStatisticsDeatiled block
$ws.Cells.Item(121, 2).Value = "Hozzáférés megtagadva"
$ws.Cells.Item(121, 3).Value = "GetStatisticsDeatiled"

$ws.Cells.Item(122, 2).Value = "Kapcsolati hiba"
$ws.Cells.Item(122, 3).Value = "GetStatisticsDeatiled"

$ws.Cells.Item(123, 2).Value = "Hiba történt az adatok mentése közben"
$ws.Cells.Item(123, 3).Value = "GetStatisticsDeatiled"

# Rows 124-126: GetTaskFilloutCount block
$ws.Cells.Item(124, 2).Value = "Hozzáférés megtagadva"
$ws.Cells.Item(124, 3).Value = "GetTaskFilloutCount"

$ws.Cells.Item(125, 3).Value = "GetTaskFilloutCount"
$ws.Cells.Item(126, 3).Value = "GetTaskFilloutCount"

# Rows 127-130: PostUserStatistic block (now 4 rows)
$ws.Cells.Item(127, 2).Value = "Hozzáférés megtagadva"
$ws.Cells.Item(127, 3).Value = "PostUserStatistic"

$ws.Cells.Item(128, 2).Value = "Kapcsolati hiba"
$ws.Cells.Item(128, 3).Value = "PostUserStatistic"

$ws.Cells.Item(129, 2).Value = "Hiba történt az adatok mentése közben"
$ws.Cells.Item(129, 3).Value = "PostUserStatistic"

$ws.Cells.Item(130, 2).Value = "Hiba történt az adatok mentése közben"
$ws.Cells.Item(130, 3).Value = "PostUserStatistic"

# Rows 131-133: GetFillingByDate block
$ws.Cells.Item(131, 2).Value = "Hozzáférés megtagadva"
$ws.Cells.Item(131, 3).Value = "GetFillingByDate"

$ws.Cells.Item(132, 3).Value = "GetFillingByDate"
$ws.Cells.Item(133, 3).Value = "GetFillingByDate"

# Rows 134, 136-138: DeleteUserStatistic block (fills in previously-empty rows 136-138)
$ws.Cells.Item(134, 2).Value = "Hozzáférés megtagadva"
$ws.Cells.Item(134, 3).Value = "DeleteUserStatistic"

$ws.Cells.Item(136, 2).Value = "Kapcsolati hiba"
$ws.Cells.Item(136, 3).Value = "DeleteUserStatistic"

$ws.Cells.Item(137, 2).Value = "Hiba történt az adatok mentése közben"
$ws.Cells.Item(137, 3).Value = "DeleteUserStatistic"

$ws.Cells.Item(138, 2).Value = "Hiba történt az adatok mentése közben"
$ws.Cells.Item(138, 3).Value = "DeleteUserStatistic"

# Update the visible window/selection to match the author's final cursor position
$ws.Activate()
$ws.Range("G130").Select()
